# Covariate_Tool dataset workbook update
# - Adds a 2nd "type of allocation" ordering: on DS1 the columns C/D/E
#   (FC-count, Elapsed-time, Failures) are rotated left by one column;
#   on DS2 the Elapsed-time/Failures columns (C/D) are swapped back.
# - DS1 becomes the active/selected sheet (was DS2), with the last used
#   cell selection moved to G11 on DS1 and G10 on DS2.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("DS1")
$ws2 = $wb.Worksheets.Item("DS2")

# --- DS1: rotate columns C,D,E left by one (new C = old D, new D = old E, new E = old C) ---
# Includes the header row (row 1) as well as all 17 data rows (rows 2-18).
for ($r = 1; $r -le 18; $r++) {
    $cVal = $ws1.Cells.Item($r, 3).Value2
    $dVal = $ws1.Cells.Item($r, 4).Value2
    $eVal = $ws1.Cells.Item($r, 5).Value2

    $ws1.Cells.Item($r, 3).Value = $dVal
    $ws1.Cells.Item($r, 4).Value = $eVal
    $ws1.Cells.Item($r, 5).Value = $cVal
}

# --- DS2: swap columns C and D (new C = old D, new D = old C, E unchanged) ---
# Includes the header row (row 1) as well as all 14 data rows (rows 2-15).
for ($r = 1; $r -le 15; $r++) {
    $cVal = $ws2.Cells.Item($r, 3).Value2
    $dVal = $ws2.Cells.Item($r, 4).Value2

    $ws2.Cells.Item($r, 3).Value = $dVal
    $ws2.Cells.Item($r, 4).Value = $cVal
}

# --- Sheet activation / selection ---
# Before: DS2 was the active (tabSelected) sheet with selection E1.
# After: DS1 is the active sheet (selection G11); DS2 keeps selection G10
# but is no longer the active tab.
$ws2.Activate()
$ws2.Range("G10").Select()

$ws1.Activate()
$ws1.Range("G11").Select()
